# DailyWorkReport.xlsx - append the 2025-01-14 work log entries (rows 30-34)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous day's block (rows 24-27, 4 rows) down onto
# the five new rows (30-34) so borders/number-formats/fonts match exactly.
$ws.Range("A24:D27").Copy()
$ws.Range("A30:D33").PasteSpecial(-4122)
$ws.Range("A27:D27").Copy()
$ws.Range("A34:D34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 30: date header + "Domm" entry
$d = Get-Date -Year 2025 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(30, 1).Value = $d.Date
$ws.Cells.Item(30, 2).Value = "Domm"
$ws.Cells.Item(30, 4).Value = 0.25

# Row 31: Study / ORM tool Demo
$ws.Cells.Item(31, 2).Value = "Study"
$ws.Cells.Item(31, 3).Value = "ORM tool Demo"
$ws.Cells.Item(31, 4).Value = 3

# Row 34 filled before row 33 to reproduce the author's shared-string order
# Row 34: Database with C# (CRUD)
$ws.Cells.Item(34, 3).Value = "Database with C# (CRUD)"
$ws.Cells.Item(34, 4).Value = 1

# Row 33: Dynamic type
$ws.Cells.Item(33, 3).Value = "Dynamic type "
$ws.Cells.Item(33, 4).Value = 2

# Row 32: Security & Cryptography
$ws.Cells.Item(32, 3).Value = "Security & Cryptography"
$ws.Cells.Item(32, 4).Value = 1.75

# Leave the selection where the author finished editing
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D33").Select()
